$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.445.30'
$ws.Range('E2').Value = '  -5.49%  '

$ws.Range('D3').Value = '2.901.59'
$ws.Range('E3').Value = '  -3.06%  '

$ws.Range('E4').Value = '  +0.17%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '548.76'
$ws.Range('E5').Value = '  -2.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '123.16'
$ws.Range('E6').Value = '  -3.77%  '

$ws.Range('E7').Value = '  +0.21%  '

$ws.Range('D8').Value = '2.889.83'
$ws.Range('E8').Value = '  -3.25%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.497'
$ws.Range('E9').Value = '  +0.34%  '

$ws.Range('E10').Value = '  -7.71%  '

$ws.Range('E11').Value = '  -9.75%  '

$ws.Range('E12').Value = '  +2.24%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000211'
$ws.Range('E13').Value = '  -4.58%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.40'
$ws.Range('E14').Value = '  -0.86%  '

$ws.Range('E15').Value = '  +1.12%  '

$ws.Range('D16').Value = '3.376.61'
$ws.Range('E16').Value = '  -3.25%  '

$ws.Range('D17').Value = '2.895.66'
$ws.Range('E17').Value = '  -3.38%  '

$ws.Range('E18').Value = '  +6.12%  '

$ws.Range('D19').Value = '57.465.68'
$ws.Range('E19').Value = '  -5.67%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '406.51'
$ws.Range('E20').Value = '  -6.98%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.93'
$ws.Range('E21').Value = '  -1.08%  '

$ws.Range('E22').Value = '  +1.48%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.83'
$ws.Range('E23').Value = '  -4.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '12.87'
$ws.Range('E24').Value = '  +3.20%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '77.10'
$ws.Range('E25').Value = '  -1.85%  '

$ws.Range('E26').Value = '  +0.06%  '

$ws.Range('E27').Value = '  +0.16%  '

$ws.Range('E28').Value = '  -1.36%  '

$ws.Range('E29').Value = '  +3.02%  '

$ws.Range('E30').Value = '  +1.11%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.04'
$ws.Range('E31').Value = '  -1.79%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '24.71'
$ws.Range('E32').Value = '  -2.87%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0974'
$ws.Range('E33').Value = '  +4.57%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.905'
$ws.Range('E34').Value = '  -3.79%  '

$ws.Range('E35').Value = '  -2.18%  '

$ws.Range('E36').Value = '  -10.74%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '47.72'
$ws.Range('E37').Value = '  -4.08%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.42'
$ws.Range('E38').Value = '  +9.27%  '

$ws.Range('D39').Value = '0.0₃0621'
$ws.Range('E39').Value = '  -6.95%  '

$ws.Range('B40').Value = 'Kaspa'
$ws.Range('C40').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.106'
$ws.Range('E40').Value = '  -1.05%  '

$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0342'
$ws.Range('E41').Value = '  -4.91%  '

$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '2.616.89'
$ws.Range('E42').Value = '  -1.08%  '

$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.42'
$ws.Range('E43').Value = '  -1.34%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '360.45'
$ws.Range('E44').Value = '  -2.45%  '

$ws.Range('E45').Value = '  +0.03%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '119.74'
$ws.Range('E46').Value = '  +1.35%  '

$ws.Range('E47').Value = '  -2.50%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.106'
$ws.Range('E48').Value = '  +0.32%  '

$ws.Range('E49').Value = '  -1.53%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.94'
$ws.Range('E50').Value = '  -1.78%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.95'
$ws.Range('E51').Value = '  -3.90%  '
